$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The old JOIN block (using generic surrogate keys "id") that appears in every
# stored SQL query, and the new JOIN block (using explicit "study_id" /
# "participant_id" keys) that should replace it everywhere it occurs.
$oldJoins = "LEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`nLEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`nLEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"`nLEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`""

$newJoins = "LEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`""

# Update every cell on the sheet whose text contains the old JOIN block,
# replacing it with the updated one (covers the StatQuery cell in column C
# and every TabQuery cell in column B).
$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $val -is [string] -and $val.Contains("df_participant prt ON std.id")) {
            $cell.Value = $val.Replace($oldJoins, $newJoins)
        }
    }
}

# Widen column C (StatQuery) and drop the explicit "best fit" flag, matching
# the updated autosize-free width used for the longer query text. 67.5
# character-units is the closest settable ColumnWidth that rounds to the
# target stored column width of ~68.332.
$ws.Columns.Item(3).ColumnWidth = 67.5
